# 自动更新Excel文件 - 2026-01-22 23:14:28
# Decrement the "剩余" (remaining) value in column E by 1 for every data row
# (rows 2..99), except row 36 which stays unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 99; $row++) {
    if ($row -eq 36) {
        continue
    }
    $cell = $ws.Cells.Item($row, 5)
    $current = $cell.Value2
    $cell.Value2 = $current - 1
}
